$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New note row (row 24): commentary about trying distance=1
$ws.Range("A24").Value = "May try also some runs with distance=1, just to see if it still works, and whether it can find peaks without double-counting."

# New data row (row 39) - results of the new run with distance set to 1
$ws.Range("A39").Value = "sg_rr_20_025 2023-12-13 17-59-26.csv"
$ws.Range("B39").Value = 0.01
$ws.Range("C39").Value = 1000
$ws.Range("D39").Value = 5001
$ws.Range("E39").Value = 1530
$ws.Range("F39").Value = 1570
$ws.Range("G39").Value = 0.5
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 5
$ws.Range("J39").Value = 3.0433333333333299
$ws.Range("K39").Value = 0.58130382883148601
$ws.Range("L39").Value = "yes but looks to have also found peaks in noise"
$ws.Range("M39").Value = "prominence setto be same as a previous run "

# K39 carries its own (non-default) font record in the target workbook -
# give it an explicit font so a second font / third cell style gets created
$ws.Range("K39").Font.Name = "Calibri"
$ws.Range("K39").Font.Size = 11

# Update the view: scroll so the new rows are visible, and select the new entry row
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("A39").Select()
